$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New bug report row (row 3)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Import loses customer type"
$ws.Range("C3").Value = "OPEN"
$ws.Range("D3").Value = "Importer"
$ws.Range("E3").Value = "The importing process does no mark what type of customer type each company is any more."
$ws.Range("F3").Value = 40245
$ws.Range("G3").Value = 40245

# Match the taller wrapped-text row height used for this entry
$ws.Rows.Item(3).RowHeight = 30

# Restore the selection to where the author left it
$null = $ws.Range("E5").Select()
